# "overhaul of whole project"
# - Adds a "test fraction" pair of columns (I/J) to the summary table,
#   computed as formulas from the existing eyes/scans + subvolumes counts.
# - Updates the Fiddler Crab row's training/testing counts (adds a
#   "testing subvolumes" count in G4, clears the old H4 value, updates
#   C4/D4/F4).
# - Moves the active selection to D11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells on row 3: "test fraction" columns ---
$ws.Range("I3").Value = "Test fraction in subvolumes/windows"
$ws.Range("J3").Value = "Test fraction in number of eyes/scans"

# --- Row 4 (Fiddler Crab) data updates ---
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 156
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 63
$ws.Range("H4").Value = ""

# --- New formula cells computing the test fractions ---
$ws.Range("I4").Formula = "=G4/(G4+D4)"
$ws.Range("J4").Formula = "=F4/(F4+C4)"

# --- Move active cell / selection to D11 ---
$ws.Range("D11").Select() | Out-Null
